$wb = $excel.ActiveWorkbook

# --- Operator sheet: update number of Type 1 / Type 2 aircraft ---
$opSheet = $wb.Worksheets.Item("Operator")
$opSheet.Range("C10").Value = 5
$opSheet.Range("D10").Value = 4

# --- Operator sheet: mark all ports as serviced (Yes) ---
$opSheet.Range("C13").Value = "Yes"
$opSheet.Range("D13").Value = "Yes"
$opSheet.Range("E13").Value = "Yes"
$opSheet.Range("F13").Value = "Yes"
$opSheet.Range("G13").Value = "Yes"

# --- Ports sheet: update Landing Slots for all ports from 1 to 3 ---
$portsSheet = $wb.Worksheets.Item("Ports")
$portsSheet.Range("E2").Value = 3
$portsSheet.Range("E3").Value = 3
$portsSheet.Range("E4").Value = 3
$portsSheet.Range("E5").Value = 3
$portsSheet.Range("E6").Value = 3

# --- Selections / active sheet, mirroring the author's final view state ---
$opSheet.Range("C4").Select()
$portsSheet.Select()
$portsSheet.Range("A1").Select()

$excel.Calculate()
